# bko2 age table.xlsx — "updating glms and adding new graphs"
#
# Semantic edit (Sheet1, the bkp id/age/range table):
#   - Row 40 (bkp061, age 6, range "6 to 8") is removed.
#   - Every row below it shifts up by one to close the gap.
#   - A brand-new record is appended at the bottom (what was row 65 is
#     now re-populated): bkp121, age 28, range "Adult".
#   - The sheet's sortState/sortCondition range shrinks from A2:C65 /
#     A2:A65 to A2:C64 / A2:A64 to match the (still-sorted) data range.
#   - The active selection moves from A10 down to A40 (whole row),
#     with the view scrolled down toward that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 41..65 up into rows 40..64 (i.e. delete old row 40's data),
# cell-by-cell, so every other row's formatting/styles stay untouched.
for ($r = 40; $r -le 64; $r++) {
    $next = $r + 1
    $ws.Range("A$r").Value2 = $ws.Range("A$next").Value2
    $ws.Range("B$r").Value2 = $ws.Range("B$next").Value2
    $ws.Range("C$r").Value2 = $ws.Range("C$next").Value2
}

# Populate the now-freed last row with the newly added record.
$ws.Range("A65").Value = "bkp121"
$ws.Range("B65").Value = 28
$ws.Range("C65").Value = "Adult"

# Re-apply the existing ID sort over the new (one-row-shorter) extent.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A64"))
$ws.Sort.SetRange($ws.Range("A2:C64"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Leave the view/selection the way the author last left it.
$ws.Activate() | Out-Null
$ws.Rows("40").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 24
